# Replace the date and each "a op b=" expression with its new value,
# walking forward through the document so duplicate old-text values
# (e.g. "24+7=" which appears twice) are each matched exactly once,
# in document order, against their corresponding new value.
$replacements = @(
    @("2025-09-20 Saturday", "2025-09-21 Sunday"),
    @("43-8=", "94-75="),
    @("70-12=", "32-6="),
    @("83-76=", "36+28="),
    @("71-12=", "73-44="),
    @("69+5=", "80-1="),
    @("84-16=", "34+58="),
    @("19+7=", "33-4="),
    @("55-38=", "38+39="),
    @("8+65=", "55+36="),
    @("43-35=", "30-21="),
    @("95-9=", "20-14="),
    @("33+59=", "27+15="),
    @("75-37=", "53-35="),
    @("35+6=", "31-14="),
    @("76+16=", "65-59="),
    @("79+8=", "77-38="),
    @("78+6=", "91-78="),
    @("24+7=", "37+45="),
    @("9+57=", "12-8="),
    @("63-15=", "17+45="),
    @("23+49=", "56+29="),
    @("48+4=", "5+28="),
    @("7+49=", "19+24="),
    @("81-55=", "9+87="),
    @("86-57=", "67-58="),
    @("5+67=", "14+67="),
    @("71-16=", "45+38="),
    @("29+69=", "96-69="),
    @("91-76=", "15+57="),
    @("5+66=", "7+68="),
    @("83-45=", "52+39="),
    @("85-37=", "33-28="),
    @("9+64=", "80-24="),
    @("82-15=", "18+43="),
    @("51-39=", "67+25="),
    @("87-68=", "6+65="),
    @("8+86=", "81-36="),
    @("93-18=", "24+69="),
    @("64-37=", "8+83="),
    @("63-55=", "71-26="),
    @("64-47=", "71-63="),
    @("90-65=", "57+4="),
    @("23+38=", "55+8="),
    @("34-29=", "71-36="),
    @("74-26=", "88+9="),
    @("78-29=", "82-13="),
    @("55+39=", "16-9="),
    @("36+8=", "82-67="),
    @("19+4=", "79+18="),
    @("78+18=", "82-28="),
    @("93-37=", "53-4="),
    @("71-54=", "22+9="),
    @("19+73=", "72-27="),
    @("55-9=", "74+18="),
    @("47+16=", "61-26="),
    @("58+5=", "59+34="),
    @("41-32=", "43-28="),
    @("72-54=", "43-18="),
    @("84-7=", "56+18="),
    @("70-19=", "30-8="),
    @("38+8=", "11-4="),
    @("6+19=", "95-78="),
    @("26+5=", "82-66="),
    @("81-3=", "27+56="),
    @("63+19=", "43-28="),
    @("28+23=", "33-18="),
    @("93-88=", "14+8="),
    @("71-39=", "59+8="),
    @("29+34=", "41-26="),
    @("81-78=", "66+6="),
    @("58-39=", "39+57="),
    @("43+19=", "95-86="),
    @("74-39=", "70-17="),
    @("24+7=", "95-79="),
    @("44+9=", "83-29="),
    @("60-38=", "92-89="),
    @("80-59=", "9+88="),
    @("47-38=", "82-74="),
    @("40-35=", "66-49="),
    @("31-6=", "19+76="),
    @("39+5=", "90-57="),
    @("26+18=", "39+44="),
    @("29+48=", "19+69="),
    @("93-79=", "27+8="),
    @("17+58=", "7+35="),
    @("53-34=", "90-59="),
    @("70-62=", "18+14="),
    @("43+28=", "91-83="),
    @("52-43=", "46+16="),
    @("93-15=", "81-68="),
    @("34-28=", "15+67="),
    @("28+55=", "62-24="),
    @("58+7=", "7+37="),
    @("44+18=", "49+45="),
    @("70-48=", "60-12="),
    @("49+12=", "84-36="),
    @("41-38=", "13+38="),
    @("29+12=", "35+39="),
    @("17+49=", "46+27="),
    @("40-34=", "17+18="),
)

$d = $word.ActiveDocument
$rng = $d.Content
$rng.Start = 0
$rng.End = $d.Content.End

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 0, $false, $newText, 1)
    if (-not $found) {
        Write-Host "MISSING: $oldText"
    }
    # Move past the just-replaced text and widen back out to the end
    # of the document so the next Find starts searching from here.
    $rng.Collapse(0)
    $rng.End = $d.Content.End
}
